$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1 "Play All Lucky Clover Slot for Free").
# --------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of All Lucky Clover online slot game and play for free with expanding Wilds and Scatters paying out regardless of position.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml) | Out-Null

# --------------------------------------------------------------------
# 2) Near the end of the document, drop the duplicated bold
#    "Play All Lucky Clover Slot for Free" paragraph entirely.
# --------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs.Item($count - 1)
$boldPara.Range.Delete() | Out-Null

# --------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    DALLE feature-image prompt, keeping the italic run formatting.
# --------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count2)
$fullRange = $lastPara.Range
$bodyRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$newText = "Feature image prompt for DALLE: Create a colorful cartoon-style image for `"All Lucky Clover`" slot game that reflects the game's fun and upbeat theme. The image should feature a happy Maya warrior wearing glasses. The warrior can be surrounded by a field of clovers or holding a clover in their hand. The illustration should be bright and cheerful, with a mix of greens, golds, and other bold colors that complement the game's overall aesthetics. The title of the game should be prominently displayed in the image, along with some of the game's symbols, such as fruits, the clover jewel, horseshoe, and diamond. The image should be eye-catching, inviting, and representative of the game's exciting features and potential payouts."
$bodyRange.Text = $newText
